# Regenerate the K (strikeout) column (column G) values for melancon_mark.xlsx
# per the updated save_data pipeline (K computed from box-score K instead of
# legacy Strike# counter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,0,0,0,0,0,0,0,1,1,1,1,0,0,0,0,1,0,1,1,1,0,0,1,2,1,1,1,1,2,2,0,2,1,0,0,0,1,0,0,0,1,1,1,1,2,1,0,0,0,0,0,1,0,1,0,0,0,2,0,0,0,0,0,2,0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}

Write-Host "Updated G2:G67 with regenerated K values."
